# Apply the commit's changes to the workbook:
#  - rename the sheet "Sheet1" -> "raw data"
#  - fix the shared string "Answer_relevance" -> "Answer relevance" (cell D1)
#  - widen columns B and C and wrap their text
#  - give the header row (row 1) a wrap+center/top style and set row heights
#    for every data row to match the new wrapped-text layout
#  - zoom the sheet to 130%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet name -----------------------------------------------------------
$ws.Name = "raw data"

# --- fix typo in header text (D1) -----------------------------------------
$ws.Range("D1").Value = "Answer relevance"

# --- column widths (B and C), matches the new <cols> block ----------------
$ws.Columns.Item(2).ColumnWidth = 50.5
$ws.Columns.Item(3).ColumnWidth = 87.16666666666667

# --- header row (row 1): keep existing center/top alignment, add wrap -----
$ws.Range("B1:C1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 16

# --- data rows (2-71): wrap text in columns B & C --------------------------
$ws.Range("B2:C71").WrapText = $true

# --- explicit row heights for rows 1-71 (Excel's computed wrap heights) ---
$rowHeights = @(16,335,256,256,288,365,304,304,365,272,365,272,380,350,80,320,365,335,380,380,304,176,350,272,288,288,365,320,80,288,256,256,380,350,256,288,335,335,350,304,350,288,80,320,240,335,380,380,256,409.6,320,335,288,288,380,272,80,320,365,272,350,365,304,304,365,304,304,272,350,288,80)

for ($i = 0; $i -lt $rowHeights.Length; $i++) {
    $ws.Rows.Item($i + 1).RowHeight = $rowHeights[$i]
}

# --- zoom the sheet view to 130% -------------------------------------------
$excel.ActiveWindow.Zoom = 130

Write-Output "edit applied"
